# Exercicio 3 - Normalizacao: add the funcionario_id / setor_id junction
# table (rows 42-48) below the existing "setor" lookup table, pushing the
# summary/answer rows further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert one blank row before the current row 48 (the
#    bordered "box-bottom" row), shifting everything from row 48 down
#    (old 48 -> 49, old 50 -> 51, ... old 56 -> 57).
$ws.Rows("48").Insert()

# 2) Borrow the look of the existing "setor" table (rows 35-40) for the
#    new table: data-row formatting from row 36 (tiled down across the
#    new rows), then the header row formatting from row 35.
$ws.Range("B36:D36").Copy()
$ws.Range("B43:D48").PasteSpecial(-4122)

$ws.Range("B35:D35").Copy()
$ws.Range("B42:D42").PasteSpecial(-4122)

# The freshly-inserted row 48 doesn't carry the "I" column border style
# that every other data row in this block has (I42:I47 = style of I41).
$ws.Range("I47").Copy()
$ws.Range("I48").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3) Header row: "funcionario_id" | "setor_id"
$ws.Range("B42").Value2 = "funcionario_id"
$ws.Range("C42").Value2 = "setor_id"

# 4) Data rows - funcionario_id / setor_id pairs.
$ws.Range("B43").Value2 = 12
$ws.Range("C43").Value2 = 1

$ws.Range("B44").Value2 = 12
$ws.Range("C44").Value2 = 5

$ws.Range("B45").Value2 = 13
$ws.Range("C45").Value2 = 2

$ws.Range("B46").Value2 = 14
$ws.Range("C46").Value2 = 3

$ws.Range("B47").Value2 = 14
$ws.Range("C47").Value2 = 5

$ws.Range("B48").Value2 = 15
$ws.Range("C48").Value2 = 4

# 5) Merge the C:D cell of each new row, same as the "setor" table above.
$ws.Range("C42:D42").Merge()
$ws.Range("C43:D43").Merge()
$ws.Range("C44:D44").Merge()
$ws.Range("C45:D45").Merge()
$ws.Range("C46:D46").Merge()
$ws.Range("C47:D47").Merge()

# 6) Row heights for the new table rows.
$ws.Range("B42:D48").RowHeight = 15

# 7) Resize/reposition the explanatory screenshot so it still ends just
#    below the (now one-row-taller) new table.
$shp = $ws.Shapes.Item(1)
$shp.Left = 446.0314960629921
$shp.Top = 514.4031496062992
$shp.Width = 501.703937007874
$shp.Height = 211.46456692913387

# 8) Leave the selection where the author left it.
$ws.Range("D48").Select()
